# Selenium reporting workbook - add a TimeStamp column and a Test Data
# (failure/assert message) column to the Chapter1 results sheet, and
# record the outcome of the run (added Timestamp func and error handling).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chapter1")

# New headers for the timestamp + test-data columns
$ws.Range("E1").Value = "TimeStamp"
$ws.Range("F1").Value = "Test Data"

# Status column now gets populated per test case
$ws.Range("D2").Value = "PASS"
$ws.Range("D3").Value = "PASS"

# TestCase1's description was updated to reflect the actual field under test
$ws.Range("B2").Value = "This test case is to test the txt field"

# Result message + timestamps captured by the run
$ws.Range("F2").Value = "Test Message from Excel"
$ws.Range("E2").Value = "2020-06-10 16:38:59.734"
$ws.Range("E3").Value = "2020-06-10 16:39:06.907"

# Size the two new columns to fit their content
$ws.Columns.Item(5).ColumnWidth = 21.498697916666668
$ws.Columns.Item(6).ColumnWidth = 21.166666666666668

# Leave the results sheet active/selected, as it was when the run finished
$ws.Activate() | Out-Null
$ws.Range("E16").Select() | Out-Null
